$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 8, pushing the existing rows 8-11
# down to rows 10-13 (Excel default shift = shift cells down, inheriting
# the formatting of the row above, same as the native UI "Insert" action).
$ws.Rows("8:9").Insert()

# --- New row 8 ---
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 45090
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = "Otros"
$ws.Range("I8").Value = 100107001
$ws.Range("J8").Value = "Caqui"
$ws.Range("K8").Value = "Mankaki"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 150
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 17533
$ws.Range("Q8").Value = "$/caja 18 kilos granel"
$ws.Range("R8").Value = "Región del Maule"
$ws.Range("S8").Value = 974
$ws.Range("T8").Value = 18

# --- New row 9 ---
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 45090
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107001
$ws.Range("J9").Value = "Caqui"
$ws.Range("K9").Value = "Mankaki"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 130
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14462
$ws.Range("Q9").Value = "$/caja 18 kilos granel"
$ws.Range("R9").Value = "Región del Maule"
$ws.Range("S9").Value = 803
$ws.Range("T9").Value = 18
